$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 238.14285
$ws.Range("I2").Value = 254
$ws.Range("K2").Value = 254
$ws.Range("M2").Value = -141

$ws.Range("H40").Value = 4072.524
$ws.Range("I40").Value = 2105.8333
$ws.Range("J40").Value = 4859.2
$ws.Range("K40").Value = 2105.8333
$ws.Range("L40").Value = 4859.2
$ws.Range("M40").Value = -1930.8333
$ws.Range("N40").Value = -5209.2

$ws.Range("H87").Value = 142900000
$ws.Range("J87").Value = 142900000
$ws.Range("L87").Value = 142900000
$ws.Range("N87").Value = -142902496

$ws.Range("H90").Value = 142900000
$ws.Range("J90").Value = 142900000
$ws.Range("L90").Value = 428700000
$ws.Range("N90").Value = -428712480

$ws.Range("H97").Value = 3508.7778
$ws.Range("J97").Value = 3759.875
$ws.Range("L97").Value = 11279.625
$ws.Range("N97").Value = -12271.625

$ws.Range("H129").Value = 1665.2858
$ws.Range("J129").Value = 2164.5
$ws.Range("L129").Value = 6493.5
$ws.Range("N129").Value = -16493.5

$ws.Range("H138").Value = 6182635.5
$ws.Range("J138").Value = 9264451
$ws.Range("L138").Value = 27793353
$ws.Range("N138").Value = -27803633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8086.1875
$ws.Range("I61").Value = 7812.9287
$ws.Range("K61").Value = 7812.9287
$ws.Range("M61").Value = -7600.9287

$ws.Range("H88").Value = 4214.857
$ws.Range("J88").Value = 4333.3335
$ws.Range("L88").Value = 4333.3335
$ws.Range("N88").Value = -5145.3335

$ws.Range("H91").Value = 4214.857
$ws.Range("J91").Value = 4333.3335
$ws.Range("L91").Value = 4333.3335
$ws.Range("N91").Value = -7141.3335

$ws.Range("H122").Value = 7409668
$ws.Range("I122").Value = 11113042
$ws.Range("J122").Value = 2920
$ws.Range("K122").Value = 33339126
$ws.Range("L122").Value = 8760
$ws.Range("M122").Value = -33336676
$ws.Range("N122").Value = -13660

$ws.Range("H136").Value = 8086.1875
$ws.Range("I136").Value = 7812.9287
$ws.Range("K136").Value = 23438.7861
$ws.Range("M136").Value = -20888.7861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10384.846
$ws.Range("I20").Value = 16745.285
$ws.Range("J20").Value = 2964.3333
$ws.Range("K20").Value = 16745.285
$ws.Range("L20").Value = 2964.3333
$ws.Range("M20").Value = -16498.285
$ws.Range("N20").Value = -3458.3333

$ws.Range("H105").Value = 6709.45
$ws.Range("I105").Value = 7304.875
$ws.Range("K105").Value = 7304.875
$ws.Range("M105").Value = -5557.875

$ws.Range("H134").Value = 2948.5925
$ws.Range("I134").Value = 2830.65
$ws.Range("K134").Value = 8491.950000000001
$ws.Range("M134").Value = -5956.950000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 512474.5
$ws.Range("J141").Value = 512474.5
$ws.Range("L141").Value = 512474.5
$ws.Range("N141").Value = -522834.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 866776
$ws.Range("I4").Value = 744546.75
$ws.Range("J4").Value = 4778111
$ws.Range("K4").Value = 2233640.25
$ws.Range("L4").Value = 14334333
$ws.Range("M4").Value = -2233528.25
$ws.Range("N4").Value = -14334557

$ws.Range("H12").Value = 125.5625
$ws.Range("I12").Value = 74.2
$ws.Range("J12").Value = 211.16667
$ws.Range("K12").Value = 222.6
$ws.Range("L12").Value = 633.50001
$ws.Range("M12").Value = -49.60000000000002
$ws.Range("N12").Value = -979.50001

$ws.Range("H16").Value = 699
$ws.Range("I16").Value = 699
$ws.Range("K16").Value = 2097
$ws.Range("M16").Value = -1924

$ws.Range("H19").Value = 2058.4
$ws.Range("J19").Value = 199
$ws.Range("L19").Value = 597
$ws.Range("N19").Value = -945

$ws.Range("H122").Value = 903.13336
$ws.Range("J122").Value = 1163.4
$ws.Range("L122").Value = 10470.6
$ws.Range("N122").Value = -15370.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 750007500
$ws.Range("J63").Value = 750007500
$ws.Range("L63").Value = 750007500
$ws.Range("N63").Value = -750008872

$ws.Range("H66").Value = 750007500
$ws.Range("J66").Value = 750007500
$ws.Range("L66").Value = 2250022500
$ws.Range("N66").Value = -2250029364

$ws.Range("H70").Value = 20852762
$ws.Range("I70").Value = 27795518
$ws.Range("J70").Value = 24497.25
$ws.Range("K70").Value = 27795518
$ws.Range("L70").Value = 24497.25
$ws.Range("M70").Value = -27795248
$ws.Range("N70").Value = -25037.25

$ws.Range("H73").Value = 20852762
$ws.Range("I73").Value = 27795518
$ws.Range("J73").Value = 24497.25
$ws.Range("K73").Value = 27795518
$ws.Range("L73").Value = 24497.25
$ws.Range("M73").Value = -27794582
$ws.Range("N73").Value = -26369.25

$ws.Range("H80").Value = 5633.2573
$ws.Range("I80").Value = 4719.4585
$ws.Range("K80").Value = 4719.4585
$ws.Range("M80").Value = -3721.4585

$ws.Range("H83").Value = 5633.2573
$ws.Range("I83").Value = 4719.4585
$ws.Range("K83").Value = 23597.2925
$ws.Range("M83").Value = -18605.2925

$ws.Range("H117").Value = 45000
$ws.Range("J117").Value = 45000
$ws.Range("L117").Value = 45000
$ws.Range("N117").Value = -51884

$ws.Range("H122").Value = 3333.6086
$ws.Range("I122").Value = 2475.5386
$ws.Range("J122").Value = 4449.1
$ws.Range("K122").Value = 7426.6158
$ws.Range("L122").Value = 13347.3
$ws.Range("M122").Value = -4976.6158
$ws.Range("N122").Value = -18247.3

$ws.Range("H132").Value = 8277.223
$ws.Range("I132").Value = 4928.143
$ws.Range("J132").Value = 19999
$ws.Range("K132").Value = 14784.429
$ws.Range("L132").Value = 59997
$ws.Range("M132").Value = -12254.429
$ws.Range("N132").Value = -65057

$ws.Range("H136").Value = 21013.777
$ws.Range("J136").Value = 21013.777
$ws.Range("L136").Value = 63041.33099999999
$ws.Range("N136").Value = -68141.33099999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6926.3687
$ws.Range("J68").Value = 7573.2666
$ws.Range("L68").Value = 7573.2666
$ws.Range("N68").Value = -9071.266599999999

$ws.Range("H71").Value = 6926.3687
$ws.Range("J71").Value = 7573.2666
$ws.Range("L71").Value = 37866.333
$ws.Range("N71").Value = -45354.333

$ws.Range("H93").Value = 4218.1934
$ws.Range("I93").Value = 1774
$ws.Range("J93").Value = 6509.625
$ws.Range("K93").Value = 1774
$ws.Range("L93").Value = 6509.625
$ws.Range("M93").Value = -526
$ws.Range("N93").Value = -9005.625

$ws.Range("H132").Value = 4699.615
$ws.Range("I132").Value = 4098.5
$ws.Range("K132").Value = 12295.5
$ws.Range("M132").Value = -9765.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 500025000
$ws.Range("J76").Value = 500025000
$ws.Range("L76").Value = 500025000
$ws.Range("N76").Value = -500025630

$ws.Range("H79").Value = 500025000
$ws.Range("J79").Value = 500025000
$ws.Range("L79").Value = 500025000
$ws.Range("N79").Value = -500027184

$ws.Range("H107").Value = 755
$ws.Range("I107").Value = 696.75
$ws.Range("J107").Value = 813.25
$ws.Range("K107").Value = 2090.25
$ws.Range("L107").Value = 2439.75
$ws.Range("M107").Value = -170.25
$ws.Range("N107").Value = -6279.75

$ws.Range("H133").Value = 83749.75
$ws.Range("J133").Value = 83749.75
$ws.Range("L133").Value = 83749.75
$ws.Range("N133").Value = -93869.75
